# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 16798
$ws1.Range("G2").Value = "已售罄"
$ws1.Range("F3").Value = 355
$ws1.Range("F4").Value = 743
$ws1.Range("F5").Value = 255
$ws1.Range("F6").Value = 721
$ws1.Range("F7").Value = 1786

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 16798
$ws4.Range("G2").Value = "已售罄"
$ws4.Range("F3").Value = 355
$ws4.Range("F4").Value = 743
$ws4.Range("F5").Value = 255
$ws4.Range("F8").Value = 721
$ws4.Range("F9").Value = 1786
